# Auto-generated edit script: updates the crypto price/volume table
# (Sheet1) to the latest scraped values, matching the Sun Oct 27
# 2024 14:30:52 UTC GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.676.80"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.491.19"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D5").Value = "586.45"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D6").Value = "176.74"
$ws.Range("E6").Value = "  +5.07%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D9").Value = "0.140"
$ws.Range("E9").Value = "  +4.68%  "
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +3.68%  "
$ws.Range("D12").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D12").Value = "4.94"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "2.949.88"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D14").Value = "25.77"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "67.462.49"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "2.423.14"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D18").Value = "11.02"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D20").Value = "351.57"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D21").Value = "4.07"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D23").Value = "70.67"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D24").Value = "4.24"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D25").Value = "1.79"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D26").Value = "9.21"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "2.616.24"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "0.0₃0911"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D30").Value = "510.76"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D31").Value = "7.85"
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D35").Value = "0.121"
$ws.Range("E35").Value = "  +7.08%  "
$ws.Range("D36").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D36").Value = "160.90"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D39").Value = "1.34"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D41").Value = "1.72"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D42").Value = "4.88"
$ws.Range("E42").Value = "  +2.54%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D43").Value = "0.330"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D44").Value = "2.45"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("D45").NumberFormat = "@"  # keep trailing/leading zeros as text
$ws.Range("D45").Value = "143.55"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("E51").Value = "  +1.73%  "
